# greco3: filter the sale table down to entrp_ptnt_id = 1001 only.
# This reproduces applying an AutoFilter on Table1's first column
# ("entrp_ptnt_id") so that only rows whose value is 1001 remain
# visible; every other data row becomes hidden (as Excel does when a
# table AutoFilter criterion excludes a row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Filter column 1 (entrp_ptnt_id) to show only rows where the value is 1001.
$lo.Range.AutoFilter(1, @("1001"))
